# Add a "Save" column (column H) to the s_vals sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell: same style/format as the other header cells (e.g. G1),
# copy its formatting over before setting the text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Data cells: plain numeric value of 1 for each existing data row.
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
